$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1, C1, D1 were formulas (=1, =2, =3); convert them into plain text string
# constants "1", "2", "3" (shared strings), not formulas/numbers.
$rng = $ws.Range("B1:D1")
$rng.NumberFormat = "@"
$ws.Range("B1").Value = "1"
$ws.Range("C1").Value = "2"
$ws.Range("D1").Value = "3"
$rng.Style = "Normal"

# B3, C3, D3 formulas updated to multiply row1 by row2
$ws.Range("B3").Formula = "=(B1 * B2)"
$ws.Range("C3").Formula = "=(C1 * C2)"
$ws.Range("D3").Formula = "=(D1 * D2)"
